$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("sheet1")

# Insert a new column before D; existing D/E columns (and their
# formatting) shift right to E/F automatically.
$ws.Columns("D:D").Insert()

# Header row: rename the (unchanged) "last charge end time" header and
# give the newly inserted column its own header text.
$ws.Range("C1").Value = "上一次充电结束时间"
$ws.Range("D1").Value = "下一次充电开始时间"

# The inserted column has no data for the existing rows - remove the
# blank cells (and the formatting copied over by Insert) so the rows
# only span A:C plus the shifted E:F values.
$ws.Range("D2:D34").Clear()

# The (hidden) AutoFilter defined name still points at the old E column
# as the last column of the filtered range; extend it to the new last
# column F.
$n = $wb.Names.Item("sheet1!_FilterDatabase")
$n.RefersTo = "=sheet1!`$A`$1:`$F`$73"

# Match the author's final selection.
$ws.Range("D5").Select() | Out-Null
